$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B8 held a lone space " " as its value -- clear it out so the cell is
# truly empty. This also drops the now-unused shared string, which is why
# the shared-string table shrinks from 13 to 12 unique entries and every
# subsequent <v> index shifts down by one.
$ws.Range("B8").ClearContents()

# Touch (and restore) the font's Bold flag so Excel re-resolves/rewrites
# B8's cell style record -- this drops the stale applyFill/applyBorder
# flags that no longer correspond to any real formatting on the cell.
$ws.Range("B8").Font.Bold = $True
$ws.Range("B8").Font.Bold = $False

# Move/extend the active selection to match the author's final cursor
# position (B13:B14, anchored at B13).
$ws.Range("B13:B14").Select()
